$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.123.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.552.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3826"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3303"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07351"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.793"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.730"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.564.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001069"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06639"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.29%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.348"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("E23").Value = "  -2.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.129.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.294"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.531"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.12"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.933"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.740.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.070"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.875"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.900"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08215"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.288"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.49%  "

$ws.Range("E37").Value = "  -1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02320"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.282"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2146"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.233"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.14%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6024"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.731"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5837"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.967"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.08%  "

$ws.Range("E50").Value = "  -3.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07015"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.01%  "
